$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values per repulled/recalculated data
$ws.Range("F3").Value = 2
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -5
$ws.Range("F11").Value = -1
